$d = $word.ActiveDocument
$d.Content.Find.Execute("Ceuta", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ceuta,", 2)
